# Fruta / hortaliza, semanal
# Shifts the Pera records in rows 569-600 down by two rows (571-602),
# gives rows 569-570 a new reporting date, and back-fills the vacated
# rows 569-570 data window accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 569
$lastRow = 600
$lastCol = 20   # column T

# 1) Snapshot the current contents (values + D-column number format) of
#    rows 569-600 before any writes happen, since the shift reads from
#    rows that will also be written to.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}
$dateFormat = $ws.Cells.Item($firstRow, 4).NumberFormat

# 2) Shift rows 571-600 so each becomes a copy of the row two above it
#    in the original data (row 600 <- old 598, ... row 571 <- old 569).
for ($r = $lastRow; $r -ge ($firstRow + 2); $r--) {
    $src = $snapshot[$r - 2]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $src[$c]
    }
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
}

# 3) The two rows that fell off the bottom of the window (old 599, 600)
#    land in two brand-new rows, 601 and 602.
$ws.Cells.Item(601, 4).NumberFormat = $dateFormat
$ws.Cells.Item(602, 4).NumberFormat = $dateFormat
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(601, $c).Value2 = $snapshot[599][$c]
    $ws.Cells.Item(602, $c).Value2 = $snapshot[600][$c]
}

# 4) Row 569 keeps its original data except for a new reporting date.
$ws.Cells.Item(569, 4).Value2 = 44714
$ws.Cells.Item(569, 4).NumberFormat = $dateFormat

# 5) Row 570 gets a new reporting date plus updated volume/origin.
$ws.Cells.Item(570, 4).Value2 = 44714
$ws.Cells.Item(570, 4).NumberFormat = $dateFormat
$ws.Cells.Item(570, 13).Value2 = 250
$ws.Cells.Item(570, 18).Value2 = "Región de O'Higgins"
